# Language/language.xlsx edit:
#  - olivineDesc (row15) and quartzDesc (row17) "Value" cells become the
#    generic placeholder "mineral" (was "It's green." / "Colorful and
#    crystally.").
#  - Eight new rows are inserted after row 17 (before the old "calcite" row)
#    to seed four new minerals: feldspar, pyroxene, amphibole, biotite -
#    each with a Key/Value/VoiceDuration triple, matching the existing
#    mineral pattern (Name row then "<name>Desc" / "mineral" row).
#  - Selection/view state is updated to match the author's last position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update the two existing "Desc" rows to point at the new generic
#     "mineral" placeholder string ---
$ws.Range("B15").Value = "mineral"
$ws.Range("B17").Value = "mineral"

# --- make room for the four new minerals: insert 8 blank rows starting at
#     row 18 (pushes calcite and everything after it down by 8) ---
$ws.Rows("18:25").Insert()

# --- feldspar ---
$ws.Range("A18").Value = "feldspar"
$ws.Range("B18").Value = "Feldspar"
$ws.Range("C18").Value = 1

$ws.Range("A19").Value = "feldsparDesc"
$ws.Range("B19").Value = "mineral"
$ws.Range("C19").Value = 5

# --- pyroxene ---
$ws.Range("A20").Value = "pyroxene"
$ws.Range("B20").Value = "Pyroxene"
$ws.Range("C20").Value = 1

$ws.Range("A21").Value = "pyroxeneDesc"
$ws.Range("B21").Value = "mineral"
$ws.Range("C21").Value = 5

# --- amphibole ---
$ws.Range("A22").Value = "amphibole"
$ws.Range("B22").Value = "Amphibole"
$ws.Range("C22").Value = 1

$ws.Range("A23").Value = "amphiboleDesc"
$ws.Range("B23").Value = "mineral"
$ws.Range("C23").Value = 5

# --- biotite ---
$ws.Range("A24").Value = "biotite"
$ws.Range("B24").Value = "Biotite"
$ws.Range("C24").Value = 0.6

$ws.Range("A25").Value = "biotiteDesc"
$ws.Range("B25").Value = "mineral"
$ws.Range("C25").Value = 5

# --- restore selection/scroll position to match the author's last saved
#     view (best-effort; topLeftCell scroll state isn't always persisted) ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C26").Select()
